# Apply scheduled-runner value updates to each sheet (per commit diff)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 25458.7
$ws.Range("I64").Value = 32084
$ws.Range("K64").Value = 32084
$ws.Range("M64").Value = -31836

$ws.Range("H67").Value = 25458.7
$ws.Range("I67").Value = 32084
$ws.Range("K67").Value = 32084
$ws.Range("M67").Value = -31226

$ws.Range("H74").Value = 7510.6665
$ws.Range("I74").Value = 7510.6665
$ws.Range("K74").Value = 7510.6665
$ws.Range("M74").Value = -6574.6665

$ws.Range("H77").Value = 7510.6665
$ws.Range("I77").Value = 7510.6665
$ws.Range("K77").Value = 37553.3325
$ws.Range("M77").Value = -32873.3325

$ws.Range("H98").Value = 2795.8333
$ws.Range("I98").Value = 2902.647
$ws.Range("K98").Value = 2902.647
$ws.Range("M98").Value = -1404.647

$ws.Range("I121").Value = 999
$ws.Range("K121").Value = 2997
$ws.Range("M121").Value = -1250

$ws.Range("H122").Value = 2795.8333
$ws.Range("I122").Value = 2902.647
$ws.Range("K122").Value = 8707.940999999999
$ws.Range("M122").Value = -6257.940999999999

$ws.Range("H137").Value = 8593.413
$ws.Range("J137").Value = 2766.6155
$ws.Range("L137").Value = 8299.8465
$ws.Range("N137").Value = -13399.8465

$ws.Range("H138").Value = 18961.254
$ws.Range("I138").Value = 1538.7046
$ws.Range("K138").Value = 4616.1138
$ws.Range("M138").Value = 523.8861999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 807.1111
$ws.Range("I5").Value = 807.1111
$ws.Range("K5").Value = 807.1111
$ws.Range("M5").Value = -695.1111

$ws.Range("H32").Value = 54900.9
$ws.Range("I32").Value = 60501
$ws.Range("J32").Value = 4500
$ws.Range("K32").Value = 60501
$ws.Range("L32").Value = 4500
$ws.Range("M32").Value = -60214
$ws.Range("N32").Value = -5074

$ws.Range("H61").Value = 11949.4
$ws.Range("I61").Value = 1356.2858
$ws.Range("K61").Value = 1356.2858
$ws.Range("M61").Value = -1144.2858

$ws.Range("H63").Value = 2664.8333
$ws.Range("I63").Value = 1997.25
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 1997.25
$ws.Range("L63").Value = 4000
$ws.Range("M63").Value = -1311.25
$ws.Range("N63").Value = -5372

$ws.Range("H66").Value = 2664.8333
$ws.Range("I66").Value = 1997.25
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 9986.25
$ws.Range("L66").Value = 20000
$ws.Range("M66").Value = -6554.25
$ws.Range("N66").Value = -26864

$ws.Range("H136").Value = 11949.4
$ws.Range("I136").Value = 1356.2858
$ws.Range("K136").Value = 4068.8574
$ws.Range("M136").Value = -1518.8574

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 807.1111
$ws.Range("I4").Value = 807.1111
$ws.Range("K4").Value = 807.1111
$ws.Range("M4").Value = -692.1111

$ws.Range("H20").Value = 36425
$ws.Range("I20").Value = 46189.43
$ws.Range("J20").Value = 2249.5
$ws.Range("K20").Value = 46189.43
$ws.Range("L20").Value = 2249.5
$ws.Range("M20").Value = -45942.43
$ws.Range("N20").Value = -2743.5

$ws.Range("H94").Value = 2115.4285
$ws.Range("I94").Value = 2115.4285
$ws.Range("K94").Value = 2115.4285
$ws.Range("M94").Value = -1664.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2943893.8
$ws.Range("I31").Value = 4763076
$ws.Range("J31").Value = 5214.3076
$ws.Range("K31").Value = 4763076
$ws.Range("L31").Value = 5214.3076
$ws.Range("M31").Value = -4762781
$ws.Range("N31").Value = -5804.3076

$ws.Range("H34").Value = 2943893.8
$ws.Range("I34").Value = 4763076
$ws.Range("J34").Value = 5214.3076
$ws.Range("K34").Value = 4763076
$ws.Range("L34").Value = 5214.3076
$ws.Range("M34").Value = -4762874
$ws.Range("N34").Value = -5618.3076

$ws.Range("H58").Value = 15337.667
$ws.Range("I58").Value = 1809.6
$ws.Range("K58").Value = 1809.6
$ws.Range("M58").Value = -1606.6

$ws.Range("H95").Value = 27777.5
$ws.Range("J95").Value = 27777.5
$ws.Range("L95").Value = 27777.5
$ws.Range("N95").Value = -33269.5

$ws.Range("H105").Value = 2885.6667
$ws.Range("I105").Value = 2781.5715
$ws.Range("K105").Value = 2781.5715
$ws.Range("M105").Value = -1034.5715

$ws.Range("H132").Value = 85101.336
$ws.Range("I132").Value = 101121.7
$ws.Range("K132").Value = 303365.1
$ws.Range("M132").Value = -300835.1

$ws.Range("H136").Value = 15337.667
$ws.Range("I136").Value = 1809.6
$ws.Range("K136").Value = 5428.799999999999
$ws.Range("M136").Value = -2878.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 37037324
$ws.Range("J26").Value = 111111224
$ws.Range("L26").Value = 333333672
$ws.Range("N26").Value = -333334248

$ws.Range("H37").Value = 41972.555
$ws.Range("J37").Value = 41972.555
$ws.Range("L37").Value = 125917.665
$ws.Range("N37").Value = -126141.665

$ws.Range("H113").Value = 520.4737
$ws.Range("I113").Value = 353
$ws.Range("K113").Value = 1059
$ws.Range("M113").Value = 1111

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 18499
$ws.Range("J92").Value = 18499
$ws.Range("L92").Value = 18499
$ws.Range("N92").Value = -22243

$ws.Range("H102").Value = 28477.174
$ws.Range("J102").Value = 2548.3333
$ws.Range("L102").Value = 2548.3333
$ws.Range("N102").Value = -5792.3333

$ws.Range("H122").Value = 3718.4
$ws.Range("I122").Value = 3446.96
$ws.Range("K122").Value = 10340.88
$ws.Range("M122").Value = -7890.880000000001

$ws.Range("H132").Value = 4804.643
$ws.Range("I132").Value = 4772.0835
$ws.Range("K132").Value = 14316.2505
$ws.Range("M132").Value = -11786.2505

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2832.4736
$ws.Range("J46").Value = 5626.75
$ws.Range("L46").Value = 5626.75
$ws.Range("N46").Value = -6002.75

$ws.Range("H68").Value = 3961.125
$ws.Range("J68").Value = 7999
$ws.Range("L68").Value = 7999
$ws.Range("N68").Value = -9497

$ws.Range("H71").Value = 3961.125
$ws.Range("J71").Value = 7999
$ws.Range("L71").Value = 39995
$ws.Range("N71").Value = -47483

$ws.Range("H122").Value = 4203.8696
$ws.Range("I122").Value = 3256.25
$ws.Range("K122").Value = 9768.75
$ws.Range("M122").Value = -7318.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 17211
$ws.Range("J81").Value = 4880
$ws.Range("L81").Value = 9760
$ws.Range("N81").Value = -11882

$ws.Range("H84").Value = 17211
$ws.Range("J84").Value = 4880
$ws.Range("L84").Value = 48800
$ws.Range("N84").Value = -59408

$ws.Range("H122").Value = 90908.8
$ws.Range("I122").Value = 148194.78
$ws.Range("K122").Value = 444584.34
$ws.Range("M122").Value = -442134.34

$ws.Range("H132").Value = 1812.5518
$ws.Range("I132").Value = 908.0952
$ws.Range("K132").Value = 2724.2856
$ws.Range("M132").Value = -194.2856000000002
